# Insert a new weekly price record at row 236 for
# "Feria Lagunitas de Puerto Montt" / Ajo, shifting the existing rows
# 236-327 down to 237-328 (dimension grows from A1:R327 to A1:R328).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 236; everything below shifts down one row.
$ws.Rows.Item(236).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A236").Value = 4
$ws.Range("B236").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C236").Value = "Los Lagos"
$ws.Range("D236").Value = 44825
$ws.Range("E236").Value = 10
$ws.Range("F236").Value = 100112003
$ws.Range("G236").Value = "Ajo"
$ws.Range("H236").Value = "Chino"
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 80
$ws.Range("K236").Value = 23000
$ws.Range("L236").Value = 23000
$ws.Range("M236").Value = 23000
$ws.Range("N236").Value = "$/caja 10 kilos"
$ws.Range("O236").Value = "China"
$ws.Range("P236").Value = 2300
$ws.Range("Q236").Value = 10
$ws.Range("R236").Value = "Hortaliza"
